$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows ("RM 232" and "SC 92") that were dropped from the
# dataset. Delete the higher-numbered row first so the lower row index
# stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Update remaining F-column (and a few D-column) values that changed as
# part of the re-run imputation/error calculations.
$ws.Range("F2").Value = 18.03
$ws.Range("F6").Value = ""
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = ""
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = ""
$ws.Range("F24").Value = ""

$ws.Range("D26").Value = -13.8
$ws.Range("D27").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("D29").Value = -13
$ws.Range("D30").Value = -13.6
$ws.Range("D31").Value = ""
$ws.Range("F31").Value = 17.18
$ws.Range("D32").Value = ""
$ws.Range("F33").Value = 17.53
